$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing data: D2 (Number of units for Type 1) and G2 (Number of units for Type 2)
$ws.Range("D2").Value = 0
$ws.Range("G2").Value = 0

# Add new project rows (enter column-by-column to match shared-string order)

# Column A: project names
$ws.Range("A4").Value = "Sleepy"
$ws.Range("A5").Value = "Dreamy"

# Column B: neighborhoods
$ws.Range("B4").Value = "Changi Village"
$ws.Range("B5").Value = "Coney Island"

# Column C: Type 1
$ws.Range("C4").Value = "2-Room"
$ws.Range("C5").Value = "2-Room"

# Column D: Number of units for Type 1
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 1337

# Column E: Selling price for Type 1
$ws.Range("E4").Value = 350000
$ws.Range("E5").Value = 350000

# Column F: Type 2
$ws.Range("F4").Value = "3-Room"
$ws.Range("F5").Value = "3-Room"

# Column G: Number of units for Type 2
$ws.Range("G4").Value = 420
$ws.Range("G5").Value = 0

# Column H: Selling price for Type 2
$ws.Range("H4").Value = 450000
$ws.Range("H5").Value = 450000

# Column I: Application opening date
$ws.Range("I4").Value = 45703
$ws.Range("I5").Value = 45703

# Column J: Application closing date
$ws.Range("J4").Value = 45736
$ws.Range("J5").Value = 45736

# Column K: Manager
$ws.Range("K4").Value = "Donkey"
$ws.Range("K5").Value = "Monkey"

# Column L: Officer Slot
$ws.Range("L4").Value = 3
$ws.Range("L5").Value = 3

# Copy date formatting from I2/J2 to I4:J5
$ws.Range("I2:J2").Copy()
$ws.Range("I4:J5").PasteSpecial(-4122)

$ws.Range("M5").Select()
